# "Add cantrals by cantons"
#
# The sheet used to have a two-row header: row 1 carried just the unit
# labels for a few columns (m3/s / MW / GWh) and row 2 carried Hiver/Ete/
# Annee labels underneath them. This rewrites it as a single header row
# with one descriptive title per column (idx, idx2, Name, Date Start,
# Date End, (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year)
# and removes the old second header row, so the 5 plant rows shift up by
# one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old secondary header row (Hiver / Ete / Annee labels under the
# unit row) - the 5 data rows move from rows 3-7 up to rows 2-6.
$ws.Rows.Item(2).Delete()

# A1:E1 had no header text before (those columns - idx/idx2/Name/dates -
# were unlabeled) and may carry stale formatting inherited from the row
# shift, so start clean there.
$ws.Range("A1:E1").ClearFormats()

# One header per column now.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Match the font used by the rest of the labeled header/body cells (Arial 9)
# for the unit-style headers.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Put the selection on the first data row, like the source workbook.
$null = $ws.Range("A2:K2").Select()

Write-Output "rebuilt header row and removed the old units/labels row"
